# Applies the cryptos-list price/volume refresh described by the commit diff.
# Numeric-looking Price (column D) values must stay TEXT (exact formatting,
# e.g. trailing zeros) -- force text via NumberFormat "@" then restore the
# "Normal" style so no stray number-format/quote-prefix style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.432.77"
$ws.Range("E2").Value = "  -6.24%  "
$ws.Range("D3").Value = "3.257.09"
$ws.Range("E3").Value = "  -6.83%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.38%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.55%  "
$ws.Range("D9").Value = "3.252.78"
$ws.Range("E9").Value = "  -6.62%  "
$ws.Range("E10").Value = "  -11.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.58"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "629.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "3.780.84"
$ws.Range("E16").Value = "  -6.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "65.349.90"
$ws.Range("E18").Value = "  -6.29%  "
$ws.Range("E19").Value = "  -3.53%  "
$ws.Range("D20").Value = "3.258.33"
$ws.Range("E20").Value = "  -6.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.14%  "
$ws.Range("E22").Value = "  -5.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "106.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.49%  "
$ws.Range("E25").Value = "  -7.59%  "
$ws.Range("E26").Value = "  -8.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.87%  "
$ws.Range("E28").Value = "  -6.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.18%  "
$ws.Range("E31").Value = "  -10.03%  "
$ws.Range("E32").Value = "  -8.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.103"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.61"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.79%  "
$ws.Range("D36").Value = "3.711.75"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "521.37"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("E39").Value = "  -6.17%  "
$ws.Range("D40").Value = "0.0₃0726"
$ws.Range("E40").Value = "  -7.96%  "
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("B42").Value = "CoreDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.28%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "32.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("E45").Value = "  -10.77%  "
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("E47").Value = "  -7.67%  "
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.23%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  +0.41%  "
